# Auto update Excel log
# Appends the latest sensor-log rows to the ALERTS and mmWave sheets,
# matching the new log entries captured at 2026-02-01 11:3x-11:3x.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALERTS sheet: append row 21 (a new FALL_DETECTED / CRITICAL alert)
# ---------------------------------------------------------------------
$alerts = $wb.Worksheets.Item("ALERTS")

$alertsRow = 21
$alerts.Cells.Item($alertsRow, 1).Value = "'2026-02-01"
$alerts.Cells.Item($alertsRow, 2).Value = "11:34:26"
$alerts.Cells.Item($alertsRow, 3).Value = "11:00"
$alerts.Cells.Item($alertsRow, 4).Value = "Living Room"
$alerts.Cells.Item($alertsRow, 5).Value = "CRITICAL"
$alerts.Cells.Item($alertsRow, 6).Value = "FALL_DETECTED"

# ---------------------------------------------------------------------
# mmWave sheet: append rows 67-74 (new PRESENCE_DETECTED / Active events)
# ---------------------------------------------------------------------
$mmwave = $wb.Worksheets.Item("mmWave")

$mmWaveRows = @(
    @("2026-02-01", "11:34:30", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:34:41", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:37:28", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:37:39", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:37:49", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:38:01", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:38:10", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:38:21", "11:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 67
for ($i = 0; $i -lt $mmWaveRows.Count; $i++) {
    $r = $startRow + $i
    $values = $mmWaveRows[$i]
    $mmwave.Cells.Item($r, 1).Value = "'" + $values[0]
    $mmwave.Cells.Item($r, 2).Value = $values[1]
    $mmwave.Cells.Item($r, 3).Value = $values[2]
    $mmwave.Cells.Item($r, 4).Value = $values[3]
    $mmwave.Cells.Item($r, 5).Value = $values[4]
    $mmwave.Cells.Item($r, 6).Value = $values[5]
}
